$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column T (rows 2-20): correlation of each variable with v_price
$ws.Range("T2").Value = 0.0287303267030508
$ws.Range("T3").Value = -0.00495230430079656
$ws.Range("T4").Value = 0.0678001612532612
$ws.Range("T5").Value = -0.0467775591202038
$ws.Range("T6").Value = -0.0391769925902026
$ws.Range("T7").Value = -0.120653422087976
$ws.Range("T8").Value = 0.11449442715349
$ws.Range("T9").Value = 0.030974312496068
$ws.Range("T10").Value = 0.0616668743602804
$ws.Range("T11").Value = -0.0260576414308473
$ws.Range("T12").Value = 0.0185889501276473
$ws.Range("T13").Value = 0.0205364503277865
$ws.Range("T14").Value = -0.0188074646708687
$ws.Range("T15").Value = 0.0179413638887479
$ws.Range("T16").Value = 0.0183417598159164
$ws.Range("T17").Value = -0.0403656367671708
$ws.Range("T18").Value = -0.0038444135005433
$ws.Range("T19").Value = -0.0156289578168552
$ws.Range("T20").Value = 0.0620360195791417

# Mirror the same values across row 21 (v_price row), symmetric matrix
$ws.Range("A21").Value = 0.0287303267030508
$ws.Range("B21").Value = -0.00495230430079656
$ws.Range("C21").Value = 0.0678001612532612
$ws.Range("D21").Value = -0.0467775591202038
$ws.Range("E21").Value = -0.0391769925902026
$ws.Range("F21").Value = -0.120653422087976
$ws.Range("G21").Value = 0.11449442715349
$ws.Range("H21").Value = 0.030974312496068
$ws.Range("I21").Value = 0.0616668743602804
$ws.Range("J21").Value = -0.0260576414308473
$ws.Range("K21").Value = 0.0185889501276473
$ws.Range("L21").Value = 0.0205364503277865
$ws.Range("M21").Value = -0.0188074646708687
$ws.Range("N21").Value = 0.0179413638887479
$ws.Range("O21").Value = 0.0183417598159164
$ws.Range("P21").Value = -0.0403656367671708
$ws.Range("Q21").Value = -0.0038444135005433
$ws.Range("R21").Value = -0.0156289578168552
$ws.Range("S21").Value = 0.0620360195791417
